# Update Nalco PDF (2025-08-26 13:06:56 UTC)
# Appends a new "SKIPPED" run-log row (row 61) to Sheet1, mirroring the
# style of the preceding data row and extending the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 61
$prevRow = $newRow - 1

# Copy the formatting (style) of the last existing data row onto the new
# row so every cell (including the blank F/H cells) keeps style index 3.
$srcRange = $ws.Range("A$prevRow`:H$prevRow")
$dstRange = $ws.Range("A$newRow`:H$newRow")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

# Populate the new row's values.
$ws.Cells.Item($newRow, 1).Value = "2025-08-26 13:06:55 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-26 18:36:55 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 7).Value = 0
